$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores numeric-looking values (e.g. "276.40", "0.0927") as plain text
# (original cells are t="inlineStr"). Force a Text number format before assigning so Excel
# does not silently coerce the string into a floating point number and lose formatting.
$priceCells = @("D2","D3","D5","D6","D9","D10","D11","D14","D16","D17","D18","D21","D23","D24","D25","D26","D28","D29","D30","D31","D32","D33","D37","D41","D42","D43","D48","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "42.861.16"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "2.239.24"
$ws.Range("E3").Value = "  -1.60%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "114.44"
$ws.Range("E5").Value = "  +2.76%  "
$ws.Range("D6").Value = "276.40"
$ws.Range("E6").Value = "  +4.70%  "
$ws.Range("E7").Value = "  -0.95%  "
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("D9").Value = "0.609"
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("D10").Value = "46.47"
$ws.Range("E10").Value = "  -0.91%  "
$ws.Range("D11").Value = "0.0927"
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("E12").Value = "  -1.41%  "
$ws.Range("E13").Value = "  -2.94%  "
$ws.Range("D14").Value = "15.28"
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").Value = "2.577.35"
$ws.Range("E16").Value = "  -1.56%  "
$ws.Range("D17").Value = "2.236.11"
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("D18").Value = "42.537.52"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").Value = "72.24"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("E22").Value = "  -3.80%  "
$ws.Range("D23").Value = "2.98"
$ws.Range("E23").Value = "  +4.75%  "
$ws.Range("D24").Value = "231.63"
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("D25").Value = "9.27"
$ws.Range("E25").Value = "  -1.69%  "
$ws.Range("D26").Value = "12.04"
$ws.Range("E26").Value = "  +6.22%  "
$ws.Range("E27").Value = "  -1.33%  "
$ws.Range("D28").Value = "40.38"
$ws.Range("E28").Value = "  -0.89%  "
$ws.Range("D29").Value = "2.24"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "3.28"
$ws.Range("E30").Value = "  -1.83%  "
$ws.Range("D31").Value = "173.58"
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("D32").Value = "21.09"
$ws.Range("E32").Value = "  -1.82%  "
$ws.Range("D33").Value = "0.0891"
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("E35").Value = "  +11.98%  "
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("D37").Value = "0.0372"
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("E38").Value = "  +0.58%  "
$ws.Range("E39").Value = "  +1.94%  "
$ws.Range("E40").Value = "  -1.41%  "
$ws.Range("D41").Value = "70.92"
$ws.Range("E41").Value = "  -7.15%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "0.233"
$ws.Range("E42").Value = "  -1.60%  "
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").Value = "13.22"
$ws.Range("E43").Value = "  -7.36%  "
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("E45").Value = "  -2.56%  "
$ws.Range("E46").Value = "  -6.83%  "
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("D48").Value = "8.44"
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("D49").Value = "0.0989"
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").Value = "100.81"
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("D51").Value = "0.642"
$ws.Range("E51").Value = "  +7.38%  "

# Clear the temporary number format again so the saved cell style matches the original
# (unstyled) cells.
foreach ($addr in $priceCells) {
    $ws.Range($addr).ClearFormats()
}
